$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.727.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.28%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.496.00'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.45%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.45'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.26%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.27'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.59%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.616'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.49%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.490.08'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.52%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.09%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.31%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.88%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.600'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.10%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.10'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.50%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000275'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.07%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '682.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.01%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.057.31'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.38%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.82'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.62%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.700.75'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.43%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.497.84'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.41%  '

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.79%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.44'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.45%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.16'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.28%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.900'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.35%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.28'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -9.48%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.53'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.88%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -5.23%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.23%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.64'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.98%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.37'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -7.89%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.87'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -7.06%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.72'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.27%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -9.28%  '

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Mantle'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.35'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -6.72%  '

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.23'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.25%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '562.61'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.78%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.59'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -15.22%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.38%  '

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.06%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '57.04'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.01%  '

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0439'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.12%  '

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.137'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.92%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.431.88'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -8.09%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.52%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '33.36'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.37%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -8.52%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.91'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.66%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.58'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -7.64%  '

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.86%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.18'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.48%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.149'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.23%  '
